$d = $word.ActiveDocument

# 1. Replace the title words "Bookdown Clavertondown Quarto AwayDay" with "Creating Accessible Online Resources"
# Restrict the search to the title paragraph (first paragraph) to avoid touching other occurrences later in the doc.
$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Bookdown", $true, $false, $false, $false, $false, $true, 1, $false, "Creating", 2)

$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Clavertondown", $true, $false, $false, $false, $false, $true, 1, $false, "Accessible", 2)

$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Quarto", $true, $false, $false, $false, $false, $true, 1, $false, "Online", 2)

$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("AwayDay", $true, $false, $false, $false, $false, $true, 1, $false, "Resources", 2)

# 2. Fix "MathJax.As" -> "MathJax. As" (insert a missing space)
$d.Content.Find.Execute("MathJax.As", $true, $false, $false, $false, $false, $true, 1, $false, "MathJax. As", 2)
